$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before the existing "europe" column (K), pushing
# the old K column (and its data/style) into L. This also auto-fixes
# dimension, cols/col widths, and row "spans" attributes.
$ws.Columns("K:K").Insert()

# New "expected spread" data for the BerkStan (metis-4), metis-16 and
# metis-64 groups, now living in the freshly inserted column K.
$ws.Range("K17").Value = 27677
$ws.Range("K18").Value = 32239.7
$ws.Range("K19").Value = 41193.1
$ws.Range("K20").Value = 45251.9
$ws.Range("K21").Value = 67066.899999999994
$ws.Range("K22").Value = 86603.3
$ws.Range("K23").Value = 109224

$ws.Range("K24").Value = 7941.29
$ws.Range("K25").Value = 7966.33
$ws.Range("K26").Value = 10201.9
$ws.Range("K27").Value = 10075.299999999999
$ws.Range("K28").Value = 13542.3
$ws.Range("K29").Value = 15679.6
$ws.Range("K30").Value = 18982.099999999999

$ws.Range("K31").Value = 7616.64
$ws.Range("K32").Value = 7574.11
$ws.Range("K33").Value = 7925.11
$ws.Range("K34").Value = 7496.6
$ws.Range("K35").Value = 7731.6
$ws.Range("K36").Value = 7571.6
$ws.Range("K37").Value = 7950.87

# New time value for the BerkStan / 100-partition row, in a new column M,
# matching the number format used by the neighbouring L37 cell (style
# index 4 / numFmtId 2 -> "0.00").
$ws.Range("M37").Value = 1121140
$ws.Range("M37").NumberFormat = $ws.Range("L37").NumberFormat

# Match the author's final selection state.
$ws.Range("K31:K37").Select()
